# Clear the question prompt that used to live in A1 ("O que é crítico na
# sua aplicação? O que gera alto impacto?"). The header row (row 3) already
# carries the real column titles (ANÁLISE / DESENHO / MECANISMO DE
# IMPLEMENTAÇÃO / JUSTIFICATIVA), so this leftover cell text is removed,
# leaving A1 blank while keeping its existing style/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guia de Engenharia v4")
$ws.Range("A1").Value = ""
